$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header cells for Wins/Losses/Ties, matching the style used by the
# other header cells in row 1 (bold, bordered, centered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record for every data row (2 through 45).
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 77   # AD = Wins
    $ws.Cells.Item($r, 31).Value = 85   # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = Ties
}
